# Applies the odds updates for row 3 and row 4 as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value  = 2
$ws.Range("H3").Value  = 3.2
$ws.Range("I3").Value  = 4.1
$ws.Range("J3").Value  = 2.75
$ws.Range("L3").Value  = 5
$ws.Range("M3").Value  = 1.1
$ws.Range("N3").Value  = 7
$ws.Range("O3").Value  = 1.5
$ws.Range("P3").Value  = 2.5
$ws.Range("Q3").Value  = 2.6
$ws.Range("R3").Value  = 1.48
$ws.Range("S3").Value  = 1.57
$ws.Range("T3").Value  = 2.25
$ws.Range("X3").Value  = 8
$ws.Range("Y3").Value  = 9.5
$ws.Range("Z3").Value  = 17
$ws.Range("AC3").Value = 6.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 19
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 12
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.25
$ws.Range("AW3").Value = 6
$ws.Range("AX3").Value = 26
$ws.Range("AZ3").Value = 101

# Row 4 updates
$ws.Range("G4").Value  = 1.91
$ws.Range("I4").Value  = 3.9
$ws.Range("J4").Value  = 2.5
$ws.Range("L4").Value  = 4.33
$ws.Range("N4").Value  = 13
$ws.Range("W4").Value  = 7.5
$ws.Range("AC4").Value = 11
$ws.Range("AJ4").Value = 13
$ws.Range("AO4").Value = 10
$ws.Range("AU4").Value = 8
$ws.Range("AW4").Value = 6
$ws.Range("AY4").Value = 29
$ws.Range("BB4").Value = 201
